# Weekly NYPD CompStat 103rd Precinct refresh: new crime data collected
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (volume/issue number and week-covering dates) ---
$ws.Range("A8").Characters(21, 2).Text = "14"
$ws.Range("C9").Characters(27, 9).Text = "3/31/2025"
$ws.Range("C9").Characters(47, 9).Text = "4/6/2025"

# --- Data table updates (rows 14-30) ---
$ws.Range("D14").Copy($ws.Range("C14"))
$ws.Range("D15").Value = 1
$ws.Range("D15").NumberFormat = $ws.Range("F15").NumberFormat
$ws.Range("E15").Value = -100
$ws.Range("E15").NumberFormat = $ws.Range("K15").NumberFormat
$ws.Range("F15").Value = 2
$ws.Range("G15").Value = 1
$ws.Range("G15").NumberFormat = $ws.Range("F15").NumberFormat
$ws.Range("H15").Value = 100
$ws.Range("H15").NumberFormat = $ws.Range("K15").NumberFormat
$ws.Range("J15").Value = 8
$ws.Range("K15").Value = 12.5
$ws.Range("N15").Value = -43.75
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = 66.666666666666
$ws.Range("F16").Value = 19
$ws.Range("G16").Value = 20
$ws.Range("H16").Value = -5
$ws.Range("I16").Value = 71
$ws.Range("J16").Value = 78
$ws.Range("K16").Value = -8.974358974358
$ws.Range("L16").Value = -15.476190476190
$ws.Range("M16").Value = -44.53125
$ws.Range("N16").Value = -82.640586797066
$ws.Range("C17").Value = 12
$ws.Range("D17").Value = 13
$ws.Range("E17").Value = -7.692307692307
$ws.Range("F17").Value = 45
$ws.Range("G17").Value = 50
$ws.Range("H17").Value = -10
$ws.Range("I17").Value = 176
$ws.Range("J17").Value = 160
$ws.Range("K17").Value = 10
$ws.Range("L17").Value = 9.316770186335
$ws.Range("M17").Value = 104.651162790698
$ws.Range("N17").Value = -22.123893805309
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 19
$ws.Range("G18").Value = 14
$ws.Range("H18").Value = 35.714285714285
$ws.Range("I18").Value = 57
$ws.Range("J18").Value = 45
$ws.Range("K18").Value = 26.666666666666
$ws.Range("L18").Value = 18.75
$ws.Range("M18").Value = -1.724137931034
$ws.Range("N18").Value = -79.642857142857
$ws.Range("C19").Value = 5
$ws.Range("D19").Value = 7
$ws.Range("E19").Value = -28.571428571428
$ws.Range("F19").Value = 33
$ws.Range("G19").Value = 59
$ws.Range("H19").Value = -44.067796610169
$ws.Range("I19").Value = 124
$ws.Range("J19").Value = 167
$ws.Range("K19").Value = -25.748502994012
$ws.Range("L19").Value = -7.462686567164
$ws.Range("M19").Value = 19.230769230769
$ws.Range("N19").Value = -55.871886120996
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = -33.333333333333
$ws.Range("F20").Value = 18
$ws.Range("G20").Value = 17
$ws.Range("H20").Value = 5.882352941176
$ws.Range("I20").Value = 50
$ws.Range("J20").Value = 59
$ws.Range("K20").Value = -15.254237288135
$ws.Range("L20").Value = -25.373134328358
$ws.Range("M20").Value = 0
$ws.Range("N20").Value = -86.111111111111
$ws.Range("C21").Value = 28
$ws.Range("D21").Value = 31
$ws.Range("E21").Value = -9.677419354838
$ws.Range("F21").Value = 137
$ws.Range("G21").Value = 162
$ws.Range("H21").Value = -15.432098765432
$ws.Range("I21").Value = 489
$ws.Range("J21").Value = 518
$ws.Range("K21").Value = -5.598455598455
$ws.Range("L21").Value = -2.589641434262
$ws.Range("M21").Value = 12.672811059907
$ws.Range("N21").Value = -68.991756499682
$ws.Range("C22").Value = 1
$ws.Range("D22").Value = 1
$ws.Range("D22").NumberFormat = $ws.Range("F15").NumberFormat
$ws.Range("E22").Value = 0
$ws.Range("E22").NumberFormat = $ws.Range("K15").NumberFormat
$ws.Range("I22").Value = 7
$ws.Range("J22").Value = 7
$ws.Range("L22").Value = 40
$ws.Range("M22").Value = 0
$ws.Range("D14").Copy($ws.Range("C23"))
$ws.Range("D23").Value = 1
$ws.Range("D23").NumberFormat = $ws.Range("F15").NumberFormat
$ws.Range("E23").Value = -100
$ws.Range("E23").NumberFormat = $ws.Range("K15").NumberFormat
$ws.Range("F23").Value = 1
$ws.Range("G23").Value = 1
$ws.Range("G23").NumberFormat = $ws.Range("F15").NumberFormat
$ws.Range("H23").Value = 0
$ws.Range("H23").NumberFormat = $ws.Range("K15").NumberFormat
$ws.Range("J23").Value = 10
$ws.Range("K23").Value = 50
$ws.Range("M23").Value = 15.384615384615
$ws.Range("C24").Value = 30
$ws.Range("D24").Value = 40
$ws.Range("E24").Value = -25
$ws.Range("F24").Value = 169
$ws.Range("G24").Value = 117
$ws.Range("H24").Value = 44.444444444444
$ws.Range("I24").Value = 501
$ws.Range("J24").Value = 441
$ws.Range("K24").Value = 13.605442176870
$ws.Range("L24").Value = 18.720379146919
$ws.Range("M24").Value = 69.830508474576
$ws.Range("C25").Value = 16
$ws.Range("D25").Value = 25
$ws.Range("E25").Value = -36
$ws.Range("F25").Value = 104
$ws.Range("G25").Value = 76
$ws.Range("H25").Value = 36.842105263157
$ws.Range("I25").Value = 282
$ws.Range("J25").Value = 255
$ws.Range("K25").Value = 10.588235294117
$ws.Range("L25").Value = 53.260869565217
$ws.Range("C26").Value = 21
$ws.Range("D26").Value = 24
$ws.Range("E26").Value = -12.5
$ws.Range("F26").Value = 51
$ws.Range("G26").Value = 76
$ws.Range("H26").Value = -32.894736842105
$ws.Range("I26").Value = 203
$ws.Range("J26").Value = 257
$ws.Range("K26").Value = -21.011673151751
$ws.Range("L26").Value = -14.705882352941
$ws.Range("M26").Value = -9.375
$ws.Range("D27").Value = 2
$ws.Range("D27").NumberFormat = $ws.Range("F15").NumberFormat
$ws.Range("E27").Value = -100
$ws.Range("E27").NumberFormat = $ws.Range("K15").NumberFormat
$ws.Range("F27").Value = 2
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 13
$ws.Range("K27").Value = -30.769230769230
$ws.Range("C28").Value = 1
$ws.Range("D28").Value = 2
$ws.Range("E28").Value = -50
$ws.Range("G28").Value = 6
$ws.Range("H28").Value = 66.666666666666
$ws.Range("I28").Value = 30
$ws.Range("J28").Value = 29
$ws.Range("K28").Value = 3.448275862068
$ws.Range("L28").Value = 87.5
$ws.Range("F29").Value = 2
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 12
$ws.Range("K29").Value = 33.333333333333
$ws.Range("L29").Value = 200
$ws.Range("M29").Value = -7.692307692307
$ws.Range("N29").Value = -61.290322580645
$ws.Range("F30").Value = 2
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 3
$ws.Range("K30").Value = -57.142857142857
$ws.Range("L30").Value = -25
$ws.Range("M30").Value = -57.142857142857
$ws.Range("N30").Value = -89.285714285714
